{"js": "// Add a new paragraph after the existing text, with the same\n// run/paragraph formatting (inherited automatically by Word when a\n// paragraph is appended at the end of the body).\nconst body = context.document.body;\nbody.insertParagraph(\n  \"Hi, thank you for understanding just support me little more.\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n", "ps1": "# Append a new paragraph after the existing content, carrying over the\n# same language/run formatting (Word inherits it automatically since the\n# new paragraph is created right after the current last paragraph).\n$d = $word.ActiveDocument\n$d.Content.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.Text = \"Hi, thank you for understanding just support me little more.\"\n"}
